$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("F2").Value = 11.87531291856449
$ws.Range("G2").Value = 11.64171663787083
$ws.Range("H2").Value = 12.10637697944409
$ws.Range("I2").Value = 1.963892990706985
$ws.Range("J2").Value = 1.942591949804927
$ws.Range("K2").Value = 1.98512201043353
$ws.Range("L2").Value = 0.1527866467839192
$ws.Range("M2").Value = 0.1511356634833285
$ws.Range("N2").Value = 0.1544333745241097

# Row 3
$ws.Range("F3").Value = 0.001953086267067525
$ws.Range("G3").Value = 0.001202108190285474
$ws.Range("H3").Value = 0.002859145694499103
$ws.Range("I3").Value = 0.001808712551014187
$ws.Range("J3").Value = 0.001105549696174063
$ws.Range("K3").Value = 0.002654615185705069
$ws.Range("L3").Value = 0.002037019747123656
$ws.Range("M3").Value = 0.001271384674632113
$ws.Range("N3").Value = 0.002957286730969938

# Row 4
$ws.Range("F4").Value = 11.87726600483155
$ws.Range("G4").Value = 11.64291874606112
$ws.Range("H4").Value = 12.10923612513859
$ws.Range("I4").Value = 1.965701703257999
$ws.Range("J4").Value = 1.943697499501101
$ws.Range("K4").Value = 1.987776625619235
$ws.Range("L4").Value = 0.1548236665310428
$ws.Range("M4").Value = 0.1524070481579606
$ws.Range("N4").Value = 0.1573906612550796
